$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row    = 12
$srcRow = 11

# --- 1) Give the new row's label cell (column A) the same look as the
#        other year-label cells above it (bold, centered, bordered). ---
$ws.Cells.Item($srcRow, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = "2021年"

# --- 2) Columns with no data for 2021: emit them as explicit blank text
#        cells (matching the sheet's existing convention where every
#        column has a cell, even when there's no value). A bare
#        Value = "" does not create a cell in this engine, so we write a
#        quote-prefixed empty string (forces a real, empty text cell)
#        and then reset the style to the sheet default so no visible
#        formatting leaks in. ---
$emptyCols = @(3, 4, 6, 7, 11, 13, 14, 18, 21, 23, 24, 26, 30, 34, 35, 38, 39, 44, 46, 50)
foreach ($col in $emptyCols) {
    $ws.Cells.Item($row, $col).Value = "'"
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# --- 3) The actual 2021 data values. ---
$values = @{
    2  = 81
    5  = 660
    8  = 1
    9  = 3
    10 = 1
    12 = 218
    15 = 685
    16 = 1248
    17 = 21
    19 = 597
    20 = 1
    22 = 265
    25 = 1
    27 = 6589
    28 = 6
    29 = 311
    31 = 1
    32 = 12
    33 = 111
    36 = 246
    37 = 46
    40 = 2
    41 = 3
    42 = 1298
    43 = 194
    45 = 8
    47 = 499
    48 = 22
    49 = 48
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$col]
}
